$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3589764260989341
$ws.Cells.Item(2, 3).Value = 0.08360667578511993
$ws.Cells.Item(2, 4).Value = 0.05061001181340075
$ws.Cells.Item(2, 6).Value = 1.349655725240893
$ws.Cells.Item(2, 7).Value = 0.002479245429810422
$ws.Cells.Item(2, 9).Value = 0.9558457619222835
$ws.Cells.Item(2, 11).Value = 0.4297665996791409
$ws.Cells.Item(2, 13).Value = 0.8804075373533919
$ws.Cells.Item(2, 14).Value = 2.111938390888881
$ws.Cells.Item(3, 2).Value = 0.3289256157109151
$ws.Cells.Item(3, 3).Value = 0.07573555820522415
$ws.Cells.Item(3, 4).Value = 0.05049673220658057
$ws.Cells.Item(3, 6).Value = 1.332582306571013
$ws.Cells.Item(3, 7).Value = 0.002482830265372861
$ws.Cells.Item(3, 9).Value = 0.951168184292726
$ws.Cells.Item(3, 11).Value = 0.3926304161402641
$ws.Cells.Item(3, 13).Value = 0.7842511887856745
$ws.Cells.Item(3, 14).Value = 2.121121224336548
$ws.Cells.Item(4, 2).Value = 0.3106449085542238
$ws.Cells.Item(4, 3).Value = 0.07094173045319963
$ws.Cells.Item(4, 4).Value = 0.050421909669641
$ws.Cells.Item(4, 6).Value = 1.322806394695263
$ws.Cells.Item(4, 7).Value = 0.002485146351551603
$ws.Cells.Item(4, 9).Value = 0.9487337199949479
$ws.Cells.Item(4, 11).Value = 0.3700324168403597
$ws.Cells.Item(4, 13).Value = 0.7257502559203601
$ws.Cells.Item(4, 14).Value = 2.127396451580395
$ws.Cells.Item(5, 2).Value = 0.3032383976656945
$ws.Cells.Item(5, 3).Value = 0.06899796044446305
$ws.Cells.Item(5, 4).Value = 0.05039010496669682
$ws.Cells.Item(5, 6).Value = 1.319000074330376
$ws.Cells.Item(5, 7).Value = 0.002486119187485514
$ws.Cells.Item(5, 9).Value = 0.9478514642141675
$ws.Cells.Item(5, 11).Value = 0.3608748152908561
$ws.Cells.Item(5, 13).Value = 0.7020400007382648
$ws.Cells.Item(5, 14).Value = 2.130113453164419
$ws.Cells.Item(6, 2).Value = 0.3020111545649229
$ws.Cells.Item(6, 3).Value = 0.0686757859737952
$ws.Cells.Item(6, 4).Value = 0.05038474478324595
$ws.Cells.Item(6, 6).Value = 1.318378743385438
$ws.Cells.Item(6, 7).Value = 0.002486282481351303
$ws.Cells.Item(6, 9).Value = 0.947711591902646
$ws.Cells.Item(6, 11).Value = 0.3593572998545937
$ws.Cells.Item(6, 13).Value = 0.6981105613342322
$ws.Cells.Item(6, 14).Value = 2.130574249849481
$ws.Cells.Item(7, 2).Value = 0.3105448473812089
$ws.Cells.Item(7, 3).Value = 0.0709154766918374
$ws.Cells.Item(7, 4).Value = 0.05042148604572638
$ws.Cells.Item(7, 6).Value = 1.322754343382272
$ws.Cells.Item(7, 7).Value = 0.002485159353893328
$ws.Cells.Item(7, 9).Value = 0.9487213772692016
$ws.Cells.Item(7, 11).Value = 0.3699087066059121
$ws.Cells.Item(7, 13).Value = 0.7254299756410916
$ws.Cells.Item(7, 14).Value = 2.127432447481837
$ws.Cells.Item(8, 2).Value = 0.3485795597422054
$ws.Cells.Item(8, 3).Value = 0.0808845600643906
$ws.Cells.Item(8, 4).Value = 0.05057205293461919
$ws.Cells.Item(8, 6).Value = 1.343621774461141
$ws.Cells.Item(8, 7).Value = 0.002480457675919709
$ws.Cells.Item(8, 9).Value = 0.9541419710223806
$ws.Cells.Item(8, 11).Value = 0.4169197549115324
$ws.Cells.Item(8, 13).Value = 0.8471370122441186
$ws.Cells.Item(8, 14).Value = 2.11497221496488
$ws.Cells.Item(9, 2).Value = 0.4245180805817768
$ws.Cells.Item(9, 3).Value = 0.1007481758850872
$ws.Cells.Item(9, 4).Value = 0.05082503836042385
$ws.Cells.Item(9, 6).Value = 1.390176188015786
$ws.Cells.Item(9, 7).Value = 0.00247214551978132
$ws.Cells.Item(9, 9).Value = 0.9682557229324118
$ws.Cells.Item(9, 11).Value = 0.5107291909346259
$ws.Cells.Item(9, 13).Value = 1.090384298347274
$ws.Cells.Item(9, 14).Value = 2.095608480631029
$ws.Cells.Item(10, 2).Value = 0.4811394471417714
$ws.Cells.Item(10, 3).Value = 0.1155410229135327
$ws.Cells.Item(10, 4).Value = 0.05098450296551604
$ws.Cells.Item(10, 6).Value = 1.427849939661243
$ws.Cells.Item(10, 7).Value = 0.002466585667826954
$ws.Cells.Item(10, 9).Value = 0.9807680028700005
$ws.Cells.Item(10, 11).Value = 0.5806525490231991
$ws.Cells.Item(10, 13).Value = 1.272328531190368
$ws.Cells.Item(10, 14).Value = 2.084497447249376
$ws.Cells.Item(11, 2).Value = 0.5070795082767745
$ws.Cells.Item(11, 3).Value = 0.1223156049172474
$ws.Cells.Item(11, 4).Value = 0.05105118189069913
$ws.Cells.Item(11, 6).Value = 1.445750437737814
$ws.Cells.Item(11, 7).Value = 0.002464173785545226
$ws.Cells.Item(11, 9).Value = 0.9869297897579017
$ws.Cells.Item(11, 11).Value = 0.6126833767688424
$ws.Cells.Item(11, 13).Value = 1.355899441115312
$ws.Cells.Item(11, 14).Value = 2.080124652546516
$ws.Cells.Item(12, 2).Value = 0.5169285920377718
$ws.Cells.Item(12, 3).Value = 0.1248875747016314
$ws.Cells.Item(12, 4).Value = 0.05107557758200443
$ws.Cells.Item(12, 6).Value = 1.452639107305714
$ws.Cells.Item(12, 7).Value = 0.002463277235750506
$ws.Cells.Item(12, 9).Value = 0.9893309806177086
$ws.Cells.Item(12, 11).Value = 0.6248447023730819
$ws.Cells.Item(12, 13).Value = 1.387669027987187
$ws.Cells.Item(12, 14).Value = 2.078567265399272
$ws.Cells.Item(13, 2).Value = 0.5148062548219343
$ws.Cells.Item(13, 3).Value = 0.1243333609903914
$ws.Cells.Item(13, 4).Value = 0.05107036169980361
$ws.Cells.Item(13, 6).Value = 1.451150602679235
$ws.Cells.Item(13, 7).Value = 0.002463469578967398
$ws.Cells.Item(13, 9).Value = 0.9888108186723485
$ws.Cells.Item(13, 11).Value = 0.6222241228702217
$ws.Cells.Item(13, 13).Value = 1.380821282044678
$ws.Cells.Item(13, 14).Value = 2.078898288367057
$ws.Cells.Item(14, 2).Value = 0.5078892750585453
$ws.Cells.Item(14, 3).Value = 0.1225270700496424
$ws.Cells.Item(14, 4).Value = 0.05105320612090836
$ws.Cells.Item(14, 6).Value = 1.446314962533307
$ws.Cells.Item(14, 7).Value = 0.002464099690304122
$ws.Cells.Item(14, 9).Value = 0.987125975582444
$ws.Cells.Item(14, 11).Value = 0.6136832562313828
$ws.Cells.Item(14, 13).Value = 1.358510633183485
$ws.Cells.Item(14, 14).Value = 2.079994548644351
$ws.Cells.Item(15, 2).Value = 0.5036558263719826
$ws.Cells.Item(15, 3).Value = 0.1214215242439707
$ws.Cells.Item(15, 4).Value = 0.05104258626638369
$ws.Cells.Item(15, 6).Value = 1.443367351414622
$ws.Cells.Item(15, 7).Value = 0.002464487833837556
$ws.Cells.Item(15, 9).Value = 0.9861028059738643
$ws.Cells.Item(15, 11).Value = 0.6084558884365094
$ws.Cells.Item(15, 13).Value = 1.344860974677346
$ws.Cells.Item(15, 14).Value = 2.080678880861427
$ws.Cells.Item(16, 2).Value = 0.47944786880268
$ws.Cells.Item(16, 3).Value = 0.1150992049761328
$ws.Cells.Item(16, 4).Value = 0.05098002641855359
$ws.Cells.Item(16, 6).Value = 1.426695485550283
$ws.Cells.Item(16, 7).Value = 0.00246674564281424
$ws.Cells.Item(16, 9).Value = 0.9803747953317909
$ws.Cells.Item(16, 11).Value = 0.5785637308684386
$ws.Cells.Item(16, 13).Value = 1.266883786078935
$ws.Cells.Item(16, 14).Value = 2.084796976763926
$ws.Cells.Item(17, 2).Value = 0.4646437767029568
$ws.Cells.Item(17, 3).Value = 0.1112323087520508
$ws.Cells.Item(17, 4).Value = 0.05094013851091184
$ws.Cells.Item(17, 6).Value = 1.41666348133802
$ws.Cells.Item(17, 7).Value = 0.002468160718373325
$ws.Cells.Item(17, 9).Value = 0.9769813981800155
$ws.Cells.Item(17, 11).Value = 0.5602827877274308
$ws.Cells.Item(17, 13).Value = 1.219258426971081
$ws.Cells.Item(17, 14).Value = 2.087498220412371
$ws.Cells.Item(18, 2).Value = 0.4561460602680825
$ws.Cells.Item(18, 3).Value = 0.1090124260182961
$ws.Cells.Item(18, 4).Value = 0.05091664491308201
$ws.Cells.Item(18, 6).Value = 1.410965079991342
$ws.Cells.Item(18, 7).Value = 0.002468985680753877
$ws.Cells.Item(18, 9).Value = 0.9750738165367636
$ws.Cells.Item(18, 11).Value = 0.5497890019261149
$ws.Cells.Item(18, 13).Value = 1.191940752640733
$ws.Cells.Item(18, 14).Value = 2.089116037622233
$ws.Cells.Item(19, 2).Value = 0.4532718385761143
$ws.Cells.Item(19, 3).Value = 0.1082615401185194
$ws.Cells.Item(19, 4).Value = 0.0509085960698048
$ws.Cells.Item(19, 6).Value = 1.409048006187646
$ws.Cells.Item(19, 7).Value = 0.002469266899325845
$ws.Cells.Item(19, 9).Value = 0.9744355267352915
$ws.Cells.Item(19, 11).Value = 0.5462395818210553
$ws.Cells.Item(19, 13).Value = 1.182704141719341
$ws.Cells.Item(19, 14).Value = 2.089674803427528
$ws.Cells.Item(20, 2).Value = 0.4662179177881569
$ws.Cells.Item(20, 3).Value = 0.1116435056613057
$ws.Cells.Item(20, 4).Value = 0.05094444176801716
$ws.Cells.Item(20, 6).Value = 1.41772397675291
$ws.Cells.Item(20, 7).Value = 0.002468008938401489
$ws.Cells.Item(20, 9).Value = 0.9773380534078839
$ws.Cells.Item(20, 11).Value = 0.5622266594417056
$ws.Cells.Item(20, 13).Value = 1.224320393630478
$ws.Cells.Item(20, 14).Value = 2.087204027471714
$ws.Cells.Item(21, 2).Value = 0.5099202511169381
$ws.Cells.Item(21, 3).Value = 0.1230574422161226
$ws.Cells.Item(21, 4).Value = 0.05105826840154393
$ws.Cells.Item(21, 6).Value = 1.447732314159666
$ws.Cells.Item(21, 7).Value = 0.002463914156707139
$ws.Cells.Item(21, 9).Value = 0.987619010939369
$ws.Cells.Item(21, 11).Value = 0.6161910489750539
$ws.Cells.Item(21, 13).Value = 1.365060409660018
$ws.Cells.Item(21, 14).Value = 2.07966987348756
$ws.Cells.Item(22, 2).Value = 0.5386345907023156
$ws.Cells.Item(22, 3).Value = 0.1305555145111441
$ws.Cells.Item(22, 4).Value = 0.05112767681145058
$ws.Cells.Item(22, 6).Value = 1.467986707984096
$ws.Cells.Item(22, 7).Value = 0.00246133572908048
$ws.Cells.Item(22, 9).Value = 0.9947338638159877
$ws.Cells.Item(22, 11).Value = 0.6516460948743941
$ws.Cells.Item(22, 13).Value = 1.457763291514866
$ws.Cells.Item(22, 14).Value = 2.075320199408083
$ws.Cells.Item(23, 2).Value = 0.5232953102449756
$ws.Cells.Item(23, 3).Value = 0.1265501124665036
$ws.Cells.Item(23, 4).Value = 0.05109109192709127
$ws.Cells.Item(23, 6).Value = 1.457117630287541
$ws.Cells.Item(23, 7).Value = 0.002462702971081501
$ws.Cells.Item(23, 9).Value = 0.9909002346154665
$ws.Cells.Item(23, 11).Value = 0.6327060402578866
$ws.Cells.Item(23, 13).Value = 1.408217489974319
$ws.Cells.Item(23, 14).Value = 2.07758899281491
$ws.Cells.Item(24, 2).Value = 0.4655062070011127
$ws.Cells.Item(24, 3).Value = 0.111457593429833
$ws.Cells.Item(24, 4).Value = 0.05094249801374318
$ws.Cells.Item(24, 6).Value = 1.417244311483486
$ws.Cells.Item(24, 7).Value = 0.002468077522495148
$ws.Cells.Item(24, 9).Value = 0.9771766746451576
$ws.Cells.Item(24, 11).Value = 0.5613477846656565
$ws.Cells.Item(24, 13).Value = 1.222031683529451
$ws.Cells.Item(24, 14).Value = 2.08733683018535
$ws.Cells.Item(25, 2).Value = 0.4038294134049636
$ws.Cells.Item(25, 3).Value = 0.0953401291125715
$ws.Cells.Item(25, 4).Value = 0.05076119074236729
$ws.Cells.Item(25, 6).Value = 1.376975202016894
$ws.Cells.Item(25, 7).Value = 0.002474297646926749
$ws.Cells.Item(25, 9).Value = 0.9640626504190948
$ws.Cells.Item(25, 11).Value = 0.4851763492342513
$ws.Cells.Item(25, 13).Value = 1.02404443596437
$ws.Cells.Item(25, 14).Value = 2.100301548883252
